# Slide 3 ("Creating a Conda Environment (yaml file)") has a text box
# (shape 3) listing the step-by-step console commands for creating /
# activating the conda environment. Two of those lines reference the old
# environment name "wrds2_env" (and the file "wrds2_env.yml"); rename them
# to "wrds_workshop_env" / "wrds_workshop_env.yml", leaving the rest of the
# line text and all run formatting untouched.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(3)
$shp = $s.Shapes.Item(3)
$tr  = $shp.TextFrame.TextRange

# --- paragraph 6: "micromamba env create -f wrds2_env.yml"
$para6 = $tr.Paragraphs(6, 1)
$text6 = $para6.Text
$localIdx6 = $text6.IndexOf("wrds2_env.yml")
if ($localIdx6 -ge 0) {
    $absStart6 = $para6.Start + $localIdx6
    $target6 = $tr.Characters($absStart6, 13)   # "wrds2_env.yml".Length
    $target6.Text = "wrds_workshop_env.yml"
}

# --- paragraph 7: "micromamba activate wrds2_env"
$para7 = $tr.Paragraphs(7, 1)
$text7 = $para7.Text
$localIdx7 = $text7.IndexOf("wrds2_env")
if ($localIdx7 -ge 0) {
    $absStart7 = $para7.Start + $localIdx7
    $target7 = $tr.Characters($absStart7, 9)    # "wrds2_env".Length
    $target7.Text = "wrds_workshop_env"
}
